$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 143
$lastCol = 10

for ($newRow = 144; $newRow -le 147; $newRow++) {
    # Copy numeric values from the source row into the new row
    for ($c = 1; $c -le $lastCol; $c++) {
        $v = $ws.Cells.Item($srcRow, $c).Value2
        $ws.Cells.Item($newRow, $c).Value2 = $v
    }

    # Increment the date serial in column A relative to the source row
    $ws.Cells.Item($newRow, 1).Value2 = $ws.Cells.Item($srcRow, 1).Value2 + ($newRow - $srcRow)

    # Copy the formatting (styles) of the source row onto the new row
    $ws.Range("A$srcRow`:J$srcRow").Copy()
    $ws.Range("A$newRow`:J$newRow").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
